$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Conditions & Exceptions")
$ws.Select()

# Copy the formatting (borders/fonts/alignment) of the last existing row (36)
# down into the new row 37, then set the row height & values.
$ws.Range("A36:K36").Copy()
$ws.Range("A37:K37").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Rows("37:37").RowHeight = 68

$ws.Range("K37").Value = "Note: this can't be created through the rule wizard. It requires the user to setup a Microsoft Office compatible alert source and to create a rule for the alert."
$ws.Range("A37").Value = "which is an <Alert E-mail>"
$ws.Range("B37").Value = "Y"
$ws.Range("C37").Value = "N"
$ws.Range("D37").Value = "N"
$ws.Range("E37").Value = "N"
$ws.Range("F37").Value = "N"
$ws.Range("G37").Value = "N"
$ws.Range("H37").Value = "N"
$ws.Range("I37").Value = "N"
$ws.Range("J37").Value = "Y"

# Re-create conditional formatting for the new row, matching the pattern used
# for row 36 (Y -> green, N -> red), with each new rule pushed to top priority
# (mirrors how Excel re-numbers rules when the formatting is duplicated down
# to a newly added row).

$blocks = @("B37:C37,E37:G37,J37", "D37", "K37", "K37", "I37", "H37")

foreach ($addr in $blocks) {
    $rng = $ws.Range($addr)

    $fcY = $rng.FormatConditions.Add(8, 3, '"Y"')
    $fcY.Font.Color = 24832
    $fcY.Interior.Color = 13561798
    $fcY.SetFirstPriority()

    $fcN = $rng.FormatConditions.Add(8, 3, '"N"')
    $fcN.Font.Color = 393372
    $fcN.Interior.Color = 13551615
    $fcN.SetFirstPriority()
}
